$d = $word.ActiveDocument
$mult = [string][char]215

# Update date paragraph
$d.Paragraphs.Item(1).Range.Text = '2023-07-29 Saturday'

$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = '82' + $mult + '56=4592'
$t.Cell(1,2).Range.Text = '45' + $mult + '52=2340'
$t.Cell(1,3).Range.Text = '94' + $mult + '66=6204'
$t.Cell(1,4).Range.Text = '38' + $mult + '42=1596'
$t.Cell(1,5).Range.Text = '39' + $mult + '18=702'
$t.Cell(2,1).Range.Text = '92' + $mult + '73=6716'
$t.Cell(2,2).Range.Text = '26' + $mult + '63=1638'
$t.Cell(2,3).Range.Text = '98' + $mult + '63=6174'
$t.Cell(2,4).Range.Text = '88' + $mult + '92=8096'
$t.Cell(2,5).Range.Text = '28' + $mult + '49=1372'
$t.Cell(3,1).Range.Text = '44' + $mult + '29=1276'
$t.Cell(3,2).Range.Text = '92' + $mult + '69=6348'
$t.Cell(3,3).Range.Text = '67' + $mult + '66=4422'
$t.Cell(3,4).Range.Text = '43' + $mult + '44=1892'
$t.Cell(3,5).Range.Text = '50' + $mult + '67=3350'
$t.Cell(4,1).Range.Text = '100' + $mult + '15=1500'
$t.Cell(4,2).Range.Text = '94' + $mult + '42=3948'
$t.Cell(4,3).Range.Text = '28' + $mult + '99=2772'
$t.Cell(4,4).Range.Text = '95' + $mult + '44=4180'
$t.Cell(4,5).Range.Text = '100' + $mult + '45=4500'
$t.Cell(5,1).Range.Text = '49' + $mult + '58=2842'
$t.Cell(5,2).Range.Text = '17' + $mult + '39=663'
$t.Cell(5,3).Range.Text = '76' + $mult + '31=2356'
$t.Cell(5,4).Range.Text = '83' + $mult + '67=5561'
$t.Cell(5,5).Range.Text = '87' + $mult + '65=5655'
$t.Cell(6,1).Range.Text = '50' + $mult + '70=3500'
$t.Cell(6,2).Range.Text = '69' + $mult + '31=2139'
$t.Cell(6,3).Range.Text = '55' + $mult + '65=3575'
$t.Cell(6,4).Range.Text = '37' + $mult + '28=1036'
$t.Cell(6,5).Range.Text = '26' + $mult + '15=390'
$t.Cell(7,1).Range.Text = '90' + $mult + '50=4500'
$t.Cell(7,2).Range.Text = '62' + $mult + '37=2294'
$t.Cell(7,3).Range.Text = '45' + $mult + '76=3420'
$t.Cell(7,4).Range.Text = '55' + $mult + '28=1540'
$t.Cell(7,5).Range.Text = '49' + $mult + '47=2303'
$t.Cell(8,1).Range.Text = '42' + $mult + '17=714'
$t.Cell(8,2).Range.Text = '11' + $mult + '86=946'
$t.Cell(8,3).Range.Text = '13' + $mult + '67=871'
$t.Cell(8,4).Range.Text = '91' + $mult + '98=8918'
$t.Cell(8,5).Range.Text = '47' + $mult + '45=2115'
$t.Cell(9,1).Range.Text = '41' + $mult + '41=1681'
$t.Cell(9,2).Range.Text = '14' + $mult + '14=196'
$t.Cell(9,3).Range.Text = '89' + $mult + '18=1602'
$t.Cell(9,4).Range.Text = '43' + $mult + '66=2838'
$t.Cell(9,5).Range.Text = '90' + $mult + '18=1620'
$t.Cell(10,1).Range.Text = '42' + $mult + '67=2814'
$t.Cell(10,2).Range.Text = '72' + $mult + '43=3096'
$t.Cell(10,3).Range.Text = '93' + $mult + '45=4185'
$t.Cell(10,4).Range.Text = '50' + $mult + '63=3150'
$t.Cell(10,5).Range.Text = '62' + $mult + '99=6138'
$t.Cell(11,1).Range.Text = '18' + $mult + '60=1080'
$t.Cell(11,2).Range.Text = '15' + $mult + '82=1230'
$t.Cell(11,3).Range.Text = '41' + $mult + '16=656'
$t.Cell(11,4).Range.Text = '15' + $mult + '91=1365'
$t.Cell(11,5).Range.Text = '73' + $mult + '24=1752'
$t.Cell(12,1).Range.Text = '72' + $mult + '35=2520'
$t.Cell(12,2).Range.Text = '20' + $mult + '38=760'
$t.Cell(12,3).Range.Text = '51' + $mult + '51=2601'
$t.Cell(12,4).Range.Text = '90' + $mult + '92=8280'
$t.Cell(12,5).Range.Text = '50' + $mult + '23=1150'
$t.Cell(13,1).Range.Text = '47' + $mult + '62=2914'
$t.Cell(13,2).Range.Text = '56' + $mult + '46=2576'
$t.Cell(13,3).Range.Text = '11' + $mult + '33=363'
$t.Cell(13,4).Range.Text = '65' + $mult + '90=5850'
$t.Cell(13,5).Range.Text = '18' + $mult + '76=1368'
$t.Cell(14,1).Range.Text = '34' + $mult + '85=2890'
$t.Cell(14,2).Range.Text = '83' + $mult + '93=7719'
$t.Cell(14,3).Range.Text = '60' + $mult + '67=4020'
$t.Cell(14,4).Range.Text = '85' + $mult + '11=935'
$t.Cell(14,5).Range.Text = '93' + $mult + '85=7905'
$t.Cell(15,1).Range.Text = '17' + $mult + '34=578'
$t.Cell(15,2).Range.Text = '26' + $mult + '72=1872'
$t.Cell(15,3).Range.Text = '45' + $mult + '53=2385'
$t.Cell(15,4).Range.Text = '89' + $mult + '45=4005'
$t.Cell(15,5).Range.Text = '66' + $mult + '31=2046'
$t.Cell(16,1).Range.Text = '94' + $mult + '77=7238'
$t.Cell(16,2).Range.Text = '15' + $mult + '79=1185'
$t.Cell(16,3).Range.Text = '32' + $mult + '28=896'
$t.Cell(16,4).Range.Text = '34' + $mult + '42=1428'
$t.Cell(16,5).Range.Text = '23' + $mult + '66=1518'
$t.Cell(17,1).Range.Text = '64' + $mult + '77=4928'
$t.Cell(17,2).Range.Text = '19' + $mult + '59=1121'
$t.Cell(17,3).Range.Text = '29' + $mult + '38=1102'
$t.Cell(17,4).Range.Text = '41' + $mult + '18=738'
$t.Cell(17,5).Range.Text = '45' + $mult + '73=3285'
$t.Cell(18,1).Range.Text = '31' + $mult + '29=899'
$t.Cell(18,2).Range.Text = '76' + $mult + '56=4256'
$t.Cell(18,3).Range.Text = '55' + $mult + '13=715'
$t.Cell(18,4).Range.Text = '51' + $mult + '65=3315'
$t.Cell(18,5).Range.Text = '68' + $mult + '19=1292'
$t.Cell(19,1).Range.Text = '83' + $mult + '92=7636'
$t.Cell(19,2).Range.Text = '64' + $mult + '47=3008'
$t.Cell(19,3).Range.Text = '28' + $mult + '47=1316'
$t.Cell(19,4).Range.Text = '96' + $mult + '71=6816'
$t.Cell(19,5).Range.Text = '57' + $mult + '43=2451'
$t.Cell(20,1).Range.Text = '61' + $mult + '78=4758'
$t.Cell(20,2).Range.Text = '53' + $mult + '27=1431'
$t.Cell(20,3).Range.Text = '59' + $mult + '76=4484'
$t.Cell(20,4).Range.Text = '66' + $mult + '31=2046'
$t.Cell(20,5).Range.Text = '31' + $mult + '82=2542'
